# Update leave module config (leave/src/main/resources/Config.xlsx)
# as part of the "updated leave module as on 3/16/2020" commit:
#   - TestCases value (B2) changes from 55 to 51
#   - Instance value (D2) changes from Automation4 to Automation5
#   - The active selection ends up on B2 (the cell that was edited)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "51"
$ws.Range("D2").Value = "Automation5"

# Leave the selection on the cell that was just edited, matching the
# resulting <selection activeCell="B2" sqref="B2"/> in the saved sheet.
$ws.Range("B2").Select()
